$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Ranking" column (AD) for rows 4-24
$ws.Range("AD4").Value = 7
$ws.Range("AD5").Value = 4
$ws.Range("AD6").Value = 4
$ws.Range("AD7").Value = 4
$ws.Range("AD8").Value = 5
$ws.Range("AD9").Value = 7
$ws.Range("AD10").Value = 4
$ws.Range("AD11").Value = 4
$ws.Range("AD12").Value = 5
$ws.Range("AD13").Value = 11
$ws.Range("AD14").Value = 13
$ws.Range("AD15").Value = 4
$ws.Range("AD16").Value = 10
$ws.Range("AD17").Value = 8
$ws.Range("AD18").Value = 14
$ws.Range("AD19").Value = 12
$ws.Range("AD20").Value = 6
$ws.Range("AD21").Value = 3
$ws.Range("AD22").Value = 3
$ws.Range("AD23").Value = 1
$ws.Range("AD24").Value = 2

# Row 25: newly-run model result (ANNC)
$ws.Range("B25").Value = "Yes"
$ws.Range("R25").Value = 0.39
$ws.Range("S25").Value = 0.4
$ws.Range("T25").Value = "under-fitting (high bias and high variance)"
$ws.Range("U25").Value = 0.61
$ws.Range("V25").Value = 0.003349
$ws.Range("W25").Value = 0.6
$ws.Range("X25").Value = 0.001403
$ws.Range("Y25").Value = 0.88
$ws.Range("Z25").Value = 0.006131
$ws.Range("AA25").Value = 7942.58
$ws.Range("AB25").Value = 0.09
$ws.Range("AC25").Value = 137427
$ws.Range("AD25").Value = 9
